# Add new "Save" column (H) with per-row 0/1 values, rows 2-72.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H, matching style of other headers (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$saveVals = @(0,0,1,1,1,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,1,0,1,0,0,0,1,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $saveVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveVals[$i]
}
